$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The TC2 block (row 20) and TC3 block (row 28) each hold one "step" of
# test-data: an action (column B) and its expected result (column D).
# This change swaps the step content between the TC2 and TC3 blocks while
# leaving the "TC2" / "TC3" id labels (B15 / B23) in place.

$tc2Action = $ws.Range("B20").Text
$tc2Expected = $ws.Range("D20").Text
$tc3Action = $ws.Range("B28").Text
$tc3Expected = $ws.Range("D28").Text

$ws.Range("B20").Value = $tc3Action
$ws.Range("D20").Value = $tc3Expected
$ws.Range("B28").Value = $tc2Action
$ws.Range("D28").Value = $tc2Expected
